$wb = $excel.ActiveWorkbook

# Drop in the final result from the RMI script directly onto the
# PDiCECpDoC sheet: replace the formula that pulled the learning-rate
# average from "Texas Notes" with the literal value itself.
$ws = $wb.Worksheets.Item("PDiCECpDoC")
[void]$ws.Activate()
[void]$ws.Range("B2").Select()
$ws.Range("B2").Value = 0.13

# The "Texas Notes" scratch-work sheet that derived that average is no
# longer needed now that its result has been dropped in as a constant.
$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("Texas Notes").Delete()

# Leave the workbook open on the "About" sheet.
[void]$wb.Worksheets.Item("About").Activate()
